$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns (in order) that receive new values: C, D, E, F, G, L, M, N
$cols = @("C", "D", "E", "F", "G", "L", "M", "N")

$data = @(
    @(0.1455765158259812, 1.002353019487146, 0.0428164257569108, 10.81725486413654, 0.002760628370459343, 0.09778339559919402, 6.314083033152457, 1.219881304098465),
    @(0.1479148473872858, 0.9812135122777477, 0.03911771959971588, 10.49191412954178, 0.002783309043547615, 0.08488838860094461, 5.959948272760073, 1.202088916787133),
    @(0.1494191025468314, 0.9694237987550594, 0.03684309937466779, 10.30487733395705, 0.002797835355867306, 0.07700459163382334, 5.747805332552161, 1.191859306147634),
    @(0.1500493478939937, 0.9649087428074949, 0.03591493544533009, 10.23174356371317, 0.002803907556520806, 0.07379893320721465, 5.662631257258937, 1.187866878109759),
    @(0.1501550420988984, 0.9641762164972363, 0.03576073036601812, 10.21978285748281, 0.002804925105925762, 0.07326702117502748, 5.648563640180669, 1.187214625912091),
    @(0.1494275323863796, 0.9693617483108881, 0.03683058727074595, 10.30387868143981, 0.002797916627428839, 0.07696133230412272, 5.746651551268855, 1.191804747365751),
    @(0.1463685666147878, 0.9948110198198492, 0.04154165276245081, 10.70236604628963, 0.002768325096300397, 0.09332928964779796, 6.190843146477903, 1.213603092400604),
    @(0.1409124632330574, 1.054632452804242, 0.05076769928289693, 11.59022447817887, 0.002714980919344092, 0.1257650460633357, 7.106655263915229, 1.261806608190341),
    @(0.1372329679869928, 1.105327492473975, 0.05756298192295262, 12.31550635248692, 0.002678531200034097, 0.1499070043788322, 7.810863789355494, 1.300481314193831),
    @(0.1356301797818098, 1.130020103315076, 0.06066381906182272, 12.66321683740307, 0.002662518599723097, 0.1609824567291582, 8.139008203685364, 1.318770638255842),
    @(0.1350334366791124, 1.139619806331098, 0.06183991660529387, 12.79761141231404, 0.002656534570359248, 0.165192027913605, 8.264474609504987, 1.32579534150986),
    @(0.135161502683065, 1.13754101480589, 0.06158653121082125, 12.76854318221928, 0.002657819833809846, 0.1642846946632659, 8.237398167339279, 1.324278063563838),
    @(0.1355808810774235, 1.130804796314919, 0.06076053673065474, 12.67421796125558, 0.002662024707772348, 0.1613284564756299, 8.149305742514741, 1.319346588423798),
    @(0.1358390900529614, 1.126711563980052, 0.06025485104134276, 12.61680098388695, 0.002664610613896634, 0.1595197648467632, 8.095506144530475, 1.316338768431905),
    @(0.1373391419164989, 1.103747820480635, 0.05736057005549355, 12.29315446097138, 0.002679588919178669, 0.1491852517659851, 7.789582721134593, 1.299299969380968),
    @(0.1382775565674841, 1.090088166343037, 0.05558782875113266, 12.09927673025652, 0.002688921741637373, 0.1428706995953348, 7.603963476625495, 1.28902464234065),
    @(0.1388239947531815, 1.082384099960564, 0.05456906964585784, 11.9894258611987, 0.002694343385663847, 0.1392473674944199, 7.49792865678728, 1.283180122216578),
    @(0.1390101583168004, 1.07980141692309, 0.05422427109834871, 11.95251283809199, 0.002696188339572838, 0.1380219890687044, 7.462149734550763, 1.281212567034288),
    @(0.1381769687788719, 1.091526363214655, 0.05577644608759869, 12.11974210465871, 0.002687922709212171, 0.1435419861035854, 7.62364697859806, 1.290111690101099),
    @(0.1354574228358771, 1.132776504045978, 0.06100309626375378, 12.70184823727936, 0.002660787492133834, 0.162196335778134, 8.17514721317599, 1.320792406060121),
    @(0.1337394766834024, 1.161195568926018, 0.06443021421868167, 13.09825298324665, 0.002643515913808909, 0.1744797024159084, 8.542652961169097, 1.341420195877191),
    @(0.1346509439896479, 1.145889134681454, 0.06259989984051373, 12.88516503952883, 0.002652692472879281, 0.1679147027848842, 8.345832086570624, 1.33035840096727),
    @(0.1382224229282478, 1.090875691240058, 0.05569117093883236, 12.11048470037446, 0.002688374197113799, 0.143238475867193, 7.614745955110095, 1.289620039793846),
    @(0.1423305725668129, 1.037317989510655, 0.04827069889948277, 11.33783549413488, 0.002728921857450903, 0.1169436332201172, 6.853707131197382, 1.24818925823314)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}
